# Apply cryptocurrency price/volume updates from the latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.725.77'
$ws.Range('E2').Value = '  +0.32%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.889.56'
$ws.Range('E3').Value = '  +0.22%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.54'
$ws.Range('E5').Value = '  +0.96%  '

$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4763'
$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2934'
$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('E9').Value = '  -0.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.04'
$ws.Range('E10').Value = '  -0.14%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07756'
$ws.Range('E11').Value = '  +0.33%  '

$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.85'
$ws.Range('E12').Value = '  -1.10%  '

$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7399'
$ws.Range('E13').Value = '  -0.70%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.887.63'
$ws.Range('E14').Value = '  +0.12%  '

$ws.Range('E15').Value = '  +1.60%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '277.44'
$ws.Range('E16').Value = '  +0.13%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.800.43'
$ws.Range('E17').Value = '  +0.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.20'
$ws.Range('E18').Value = '  -2.97%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007578'
$ws.Range('E19').Value = '  -0.44%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9998'
$ws.Range('E20').Value = '  -0.03%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.136.31'
$ws.Range('E21').Value = '  +0.67%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.334'
$ws.Range('E22').Value = '  +0.72%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9999'
$ws.Range('E23').Value = '  +0.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.247'
$ws.Range('E24').Value = '  +0.32%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.252'
$ws.Range('E25').Value = '  -1.11%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.47'
$ws.Range('E26').Value = '  +0.47%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.85'
$ws.Range('E27').Value = '  -0.89%  '

$ws.Range('E28').Value = '  -1.46%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.346'
$ws.Range('E29').Value = '  -2.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09737'
$ws.Range('E30').Value = '  -2.88%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.502'
$ws.Range('E31').Value = '  -1.22%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.312'
$ws.Range('E32').Value = '  -0.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.208'
$ws.Range('E33').Value = '  +2.72%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04887'
$ws.Range('E34').Value = '  +1.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.127'
$ws.Range('E35').Value = '  -0.52%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7005'
$ws.Range('E36').Value = '  -0.70%  '

$ws.Range('E37').Value = '  +0.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01917'
$ws.Range('E38').Value = '  +1.95%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.804'
$ws.Range('E39').Value = '  +2.43%  '

$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '76.49'
$ws.Range('E40').Value = '  +6.82%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.350'
$ws.Range('E41').Value = '  -0.20%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.032'
$ws.Range('E42').Value = '  +2.71%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4263'
$ws.Range('E43').Value = '  +0.40%  '

$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.05'
$ws.Range('E46').Value = '  -0.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.427'
$ws.Range('E47').Value = '  +1.32%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.066'
$ws.Range('E48').Value = '  -0.89%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.76'
$ws.Range('E49').Value = '  +0.15%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '923.33'

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05768'
$ws.Range('E51').Value = '  +2.20%  '
